# Actualización automática de tasas-transfi.xlsx
# Updates the "Conversión del día" note on Hoja1 and the rate figures on "tasas".

$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the Binance conversion note with the new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.44 = 30399.03 pesos`n✅ 30399.03 pesos = 7.4 = 959.54 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- tasas!N10, O10, N12, O12: update the updated exchange-rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 134.4
$wsTasas.Range("O10").Value = 4085.63
$wsTasas.Range("N12").Value = 4109
$wsTasas.Range("O12").Value = 129.7
